# Adds a "2022-Q1" fund-holdings sheet (reusing the current "总计" sheet
# object so it keeps sheetId=3) and rebuilds a fresh "总计" sheet
# (sheetId=4) that now also lists the 2022-Q1 totals row.

$xlPasteAll = -4104
$xlPasteFormats = -4122
$xlPasteValues = -4163

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)   # 2021-Q3 (style source s=1, not used)
$sheet2 = $wb.Worksheets.Item(2)   # 2021-Q4 (style source s=2)

# ---------------------------------------------------------------------
# Helper: write a value into a single cell as TEXT, even when the text
# looks numeric (e.g. "0.0070", "007872"), without allocating a brand
# new style record in styles.xml. Plain `Range.Value = "0.0070"` gets
# silently coerced to a number by the COM layer, so instead we stage
# the literal through a formula (which evaluates to a string) in a
# scratch cell, then copy/paste-special *values only* into the real
# destination - that preserves the string type cleanly.
# ---------------------------------------------------------------------
$scratch = $sheet1.Range("ZZ1")
function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $range.PasteSpecial($xlPasteValues)
    $scratch.Clear()
}

# =======================================================================
# 1. Duplicate the existing "总计" sheet (item 3) so the copy inherits
#    sheetPr/pageMargins etc. The original keeps sheetId=3 and becomes
#    "2022-Q1"; the duplicate gets the next sheetId (4) and becomes the
#    new "总计".
# =======================================================================
$fundSheet = $wb.Worksheets.Item(3)
$fundSheet.Copy($null, $fundSheet)
$totalSheet = $wb.Worksheets.Item(4)

$fundSheet.Name = "2022-Q1"
$totalSheet.Name = "总计"

# =======================================================================
# 2. "总计" sheet: shift existing data rows down by one and insert the
#    new 2022-Q1 row at the top (row 2).
# =======================================================================

# 2021-Q3 row (currently row 3) -> row 4
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B4").PasteSpecial($xlPasteAll)

# 2021-Q4 row (currently row 2) -> row 3
$totalSheet.Range("B2:D2").Copy()
$totalSheet.Range("B3").PasteSpecial($xlPasteAll)

# Re-stamp column A (row index, 0-based) with correct style + values
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A2:A4").PasteSpecial($xlPasteFormats)
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2

# New 2022-Q1 totals row
Set-TextValue $totalSheet.Range("B2") "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 1.25

# =======================================================================
# 3. "2022-Q1" sheet: clear the inherited 总计 content and rebuild it as
#    a fund-holdings table (same layout as the 2021-Q3 / 2021-Q4 sheets).
# =======================================================================
$fundSheet.Cells.Clear()

# Header row, formatted like sheet2's header (style s=2)
$sheet2.Range("B1").Copy()
$fundSheet.Range("B1:H1").PasteSpecial($xlPasteFormats)
$fundSheet.Range("B1").Value = "基金代码"
$fundSheet.Range("C1").Value = "基金名称"
$fundSheet.Range("D1").Value = "基金规模"
$fundSheet.Range("E1").Value = "股票总仓位"
$fundSheet.Range("F1").Value = "仓位占比"
$fundSheet.Range("G1").Value = "持有市值(亿元)"
$fundSheet.Range("H1").Value = "仓位排名"

# Column A (row index), formatted like sheet2's column A (style s=2)
$sheet2.Range("A2").Copy()
$fundSheet.Range("A2:A4").PasteSpecial($xlPasteFormats)
$fundSheet.Range("A2").Value = 0
$fundSheet.Range("A3").Value = 1
$fundSheet.Range("A4").Value = 2

# Row 2 - 007872 金信稳健策略灵活配置混合
Set-TextValue $fundSheet.Range("B2") "007872"
Set-TextValue $fundSheet.Range("C2") "金信稳健策略灵活配置混合"
Set-TextValue $fundSheet.Range("D2") "25.57"
Set-TextValue $fundSheet.Range("E2") "93.73"
Set-TextValue $fundSheet.Range("F2") "4.42"
Set-TextValue $fundSheet.Range("G2") "1.1302"
$fundSheet.Range("H2").Value = 9

# Row 3 - 002256 金信行业优选灵活配置混合
Set-TextValue $fundSheet.Range("B3") "002256"
Set-TextValue $fundSheet.Range("C3") "金信行业优选灵活配置混合"
Set-TextValue $fundSheet.Range("D3") "2.43"
Set-TextValue $fundSheet.Range("E3") "93.89"
Set-TextValue $fundSheet.Range("F3") "4.60"
Set-TextValue $fundSheet.Range("G3") "0.1118"
$fundSheet.Range("H3").Value = 10

# Row 4 - 002862 金信量化精选灵活配置混合
Set-TextValue $fundSheet.Range("B4") "002862"
Set-TextValue $fundSheet.Range("C4") "金信量化精选灵活配置混合"
Set-TextValue $fundSheet.Range("D4") "0.16"
Set-TextValue $fundSheet.Range("E4") "94.28"
Set-TextValue $fundSheet.Range("F4") "4.35"
Set-TextValue $fundSheet.Range("G4") "0.0070"
$fundSheet.Range("H4").Value = 6

# Restore the originally-active tab (2021-Q3)
$sheet1.Activate()
